# Update rows 228-278 on the active worksheet:
#  - SpecCode (E) changes from SOLEVUL to ARNOLAT
#  - Weight (I) becomes -1 (unknown/not recorded)
#  - Sex (J) becomes "I" (indeterminate)
#  - MatStage (K) is cleared (no maturity stage applicable)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 228; $row -le 278; $row++) {
    $ws.Cells.Item($row, 5).Value = "ARNOLAT"   # Column E - SpecCode
    $ws.Cells.Item($row, 9).Value = -1          # Column I - W(g)
    $ws.Cells.Item($row, 10).Value = "I"        # Column J - Sex
    $ws.Cells.Item($row, 11).ClearContents()    # Column K - MatStage
}
